# Apply scheduled-runner price/profit updates across the Leve tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2250
$ws.Range("I2").Value = 2866.6667
$ws.Range("K2").Value = 2866.6667
$ws.Range("M2").Value = -2753.6667
$ws.Range("H18").Value = 356.2857
$ws.Range("I18").Value = 349
$ws.Range("K18").Value = 349
$ws.Range("M18").Value = -65
$ws.Range("H98").Value = 43500664
$ws.Range("I98").Value = 50025656
$ws.Range("K98").Value = 50025656
$ws.Range("M98").Value = -50024158
$ws.Range("H116").Value = 7761.0835
$ws.Range("I116").Value = 6599.6665
$ws.Range("K116").Value = 6599.6665
$ws.Range("M116").Value = -3157.6665
$ws.Range("H122").Value = 43500664
$ws.Range("I122").Value = 50025656
$ws.Range("K122").Value = 150076968
$ws.Range("M122").Value = -150074518
$ws.Range("H125").Value = 1578
$ws.Range("J125").Value = 1578
$ws.Range("L125").Value = 14202
$ws.Range("N125").Value = -19122
$ws.Range("H132").Value = 1524.0278
$ws.Range("I132").Value = 1464.4062
$ws.Range("J132").Value = 2001
$ws.Range("K132").Value = 4393.2186
$ws.Range("L132").Value = 6003
$ws.Range("M132").Value = -1863.2186
$ws.Range("N132").Value = -11063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7264534
$ws.Range("I61").Value = 8199747
$ws.Range("J61").Value = 133533.88
$ws.Range("K61").Value = 8199747
$ws.Range("L61").Value = 133533.88
$ws.Range("M61").Value = -8199535
$ws.Range("N61").Value = -133957.88
$ws.Range("H132").Value = 4432.946
$ws.Range("I132").Value = 1979.129
$ws.Range("K132").Value = 5937.387
$ws.Range("M132").Value = -3407.387
$ws.Range("H136").Value = 7264534
$ws.Range("I136").Value = 8199747
$ws.Range("J136").Value = 133533.88
$ws.Range("K136").Value = 24599241
$ws.Range("L136").Value = 400601.64
$ws.Range("M136").Value = -24596691
$ws.Range("N136").Value = -405701.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 589804.44
$ws.Range("I31").Value = 8185.1875
$ws.Range("J31").Value = 2916281.5
$ws.Range("K31").Value = 8185.1875
$ws.Range("L31").Value = 2916281.5
$ws.Range("M31").Value = -7890.1875
$ws.Range("N31").Value = -2916871.5
$ws.Range("H34").Value = 589804.44
$ws.Range("I34").Value = 8185.1875
$ws.Range("J34").Value = 2916281.5
$ws.Range("K34").Value = 8185.1875
$ws.Range("L34").Value = 2916281.5
$ws.Range("M34").Value = -7983.1875
$ws.Range("N34").Value = -2916685.5
$ws.Range("H58").Value = 1911.2142
$ws.Range("I58").Value = 1883.2727
$ws.Range("K58").Value = 1883.2727
$ws.Range("M58").Value = -1680.2727
$ws.Range("H86").Value = 6027.5
$ws.Range("I86").Value = 6027.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6027.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4904.5
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 6027.5
$ws.Range("I89").Value = 6027.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 30137.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -24521.5
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 1348.625
$ws.Range("I122").Value = 1299
$ws.Range("K122").Value = 3897
$ws.Range("M122").Value = -1447
$ws.Range("H134").Value = 775778.6
$ws.Range("I134").Value = 1667516.4
$ws.Range("J134").Value = 11432
$ws.Range("K134").Value = 5002549.199999999
$ws.Range("L134").Value = 34296
$ws.Range("M134").Value = -5000014.199999999
$ws.Range("N134").Value = -39366
$ws.Range("H136").Value = 1911.2142
$ws.Range("I136").Value = 1883.2727
$ws.Range("K136").Value = 5649.8181
$ws.Range("M136").Value = -3099.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 513.5455
$ws.Range("I107").Value = 403.1
$ws.Range("J107").Value = 605.5833
$ws.Range("K107").Value = 1209.3
$ws.Range("L107").Value = 1816.7499
$ws.Range("M107").Value = 710.6999999999998
$ws.Range("N107").Value = -5656.7499
$ws.Range("H127").Value = 2997.5
$ws.Range("J127").Value = 2997.5
$ws.Range("L127").Value = 8992.5
$ws.Range("N127").Value = -18912.5
$ws.Range("H138").Value = 2912.375
$ws.Range("I138").Value = 1999.75
$ws.Range("K138").Value = 5999.25
$ws.Range("M138").Value = -859.25
$ws.Range("H139").Value = 2207
$ws.Range("I139").Value = 1247
$ws.Range("K139").Value = 3741
$ws.Range("M139").Value = 1399
$ws.Range("H141").Value = 196336.75
$ws.Range("I141").Value = 502915
$ws.Range("K141").Value = 1508745
$ws.Range("M141").Value = -1503565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25130728
$ws.Range("I7").Value = 100002500
$ws.Range("J7").Value = 173470.83
$ws.Range("K7").Value = 100002500
$ws.Range("L7").Value = 173470.83
$ws.Range("M7").Value = -100002388
$ws.Range("N7").Value = -173694.83
$ws.Range("H16").Value = 4709.2856
$ws.Range("I16").Value = 4797.2
$ws.Range("J16").Value = 4489.5
$ws.Range("K16").Value = 4797.2
$ws.Range("L16").Value = 4489.5
$ws.Range("M16").Value = -4627.2
$ws.Range("N16").Value = -4829.5
$ws.Range("H122").Value = 5579.3
$ws.Range("J122").Value = 7981
$ws.Range("L122").Value = 23943
$ws.Range("N122").Value = -28843
$ws.Range("H126").Value = 25130728
$ws.Range("I126").Value = 100002500
$ws.Range("J126").Value = 173470.83
$ws.Range("K126").Value = 300007500
$ws.Range("L126").Value = 520412.49
$ws.Range("M126").Value = -300005030
$ws.Range("N126").Value = -525352.49
$ws.Range("H127").Value = 115049.8
$ws.Range("J127").Value = 115049.8
$ws.Range("L127").Value = 115049.8
$ws.Range("N127").Value = -124969.8
$ws.Range("H128").Value = 152000
$ws.Range("J128").Value = 152000
$ws.Range("L128").Value = 152000
$ws.Range("N128").Value = -161960
$ws.Range("H129").Value = 74995
$ws.Range("J129").Value = 74995
$ws.Range("L129").Value = 74995
$ws.Range("N129").Value = -84995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7560.1
$ws.Range("I122").Value = 3694.2144
$ws.Range("K122").Value = 11082.6432
$ws.Range("M122").Value = -8632.643199999999
$ws.Range("H132").Value = 1524.3636
$ws.Range("I132").Value = 1362.963
$ws.Range("K132").Value = 4088.889
$ws.Range("M132").Value = -1558.889
